# Update all metric values (rows 2-26, columns B-Q) to the new
# values produced after retraining ("atualizado todo o treinamento
# para o novo lm"). Every data row shares the same new values per
# column, so we can set each column range in one shot.
# Values are parsed via [double] cast from strings to safely handle
# scientific notation literals.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    "B" = [double]"0.9999824846980448"
    "C" = [double]"0.9990072066285755"
    "D" = [double]"0.9999999999999565"
    "E" = [double]"0.999999915790956"
    "F" = [double]"0.9999999550217807"
    "G" = [double]"1.634976554226276e-05"
    "H" = [double]"0.0009267290336323482"
    "I" = [double]"2.77339343946781e-14"
    "J" = [double]"7.975920484154391e-08"
    "K" = [double]"3.987961628773915e-08"
    "L" = [double]"0.0002547310818618769"
    "M" = [double]"0.004043484331892824"
    "N" = [double]"1.000016814689877"
    "O" = [double]"0.004215623831067973"
    "P" = [double]"120.0425940012932"
    "Q" = [double]"179.767509419835"
}

foreach ($col in $newValues.Keys) {
    $value = $newValues[$col]
    $range = $ws.Range("$col" + "2:" + "$col" + "26")
    $range.Value = $value
}
